$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $refRow, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $ws.Cells.Item($refRow, $col).Style
}

# ---- Sheet: Proximity ----
$ws = $wb.Worksheets.Item("Proximity")
Set-TextCell $ws 6 1 5 '2026-02-01'
$ws.Cells.Item(6, 2).Value = '21:05:02'
$ws.Cells.Item(6, 3).Value = '21:00'
$ws.Cells.Item(6, 4).Value = 'Bedroom Door'
$ws.Cells.Item(6, 5).Value = 'EXIT'
$ws.Cells.Item(6, 6).Value = 'User EXITED Bedroom'

# ---- Sheet: mmWave(BR) ----
$ws = $wb.Worksheets.Item("mmWave(BR)")
Set-TextCell $ws 62 1 61 '2026-02-01'
$ws.Cells.Item(62, 2).Value = '21:04:23'
$ws.Cells.Item(62, 3).Value = '21:00'
$ws.Cells.Item(62, 4).Value = 'Bedroom'
$ws.Cells.Item(62, 5).Value = 9
$ws.Cells.Item(62, 6).Value = 'Occupied'
Set-TextCell $ws 63 1 61 '2026-02-01'
$ws.Cells.Item(63, 2).Value = '21:04:24'
$ws.Cells.Item(63, 3).Value = '21:00'
$ws.Cells.Item(63, 4).Value = 'Bedroom'
$ws.Cells.Item(63, 5).Value = 7
$ws.Cells.Item(63, 6).Value = 'Occupied'
Set-TextCell $ws 64 1 61 '2026-02-01'
$ws.Cells.Item(64, 2).Value = '21:04:25'
$ws.Cells.Item(64, 3).Value = '21:00'
$ws.Cells.Item(64, 4).Value = 'Bedroom'
$ws.Cells.Item(64, 5).Value = 2
$ws.Cells.Item(64, 6).Value = 'Occupied'
Set-TextCell $ws 65 1 61 '2026-02-01'
$ws.Cells.Item(65, 2).Value = '21:04:25'
$ws.Cells.Item(65, 3).Value = '21:00'
$ws.Cells.Item(65, 4).Value = 'Bedroom'
$ws.Cells.Item(65, 5).Value = 4
$ws.Cells.Item(65, 6).Value = 'Occupied'
Set-TextCell $ws 66 1 61 '2026-02-01'
$ws.Cells.Item(66, 2).Value = '21:04:26'
$ws.Cells.Item(66, 3).Value = '21:00'
$ws.Cells.Item(66, 4).Value = 'Bedroom'
$ws.Cells.Item(66, 5).Value = 2
$ws.Cells.Item(66, 6).Value = 'Occupied'
Set-TextCell $ws 67 1 61 '2026-02-01'
$ws.Cells.Item(67, 2).Value = '21:04:28'
$ws.Cells.Item(67, 3).Value = '21:00'
$ws.Cells.Item(67, 4).Value = 'Bedroom'
$ws.Cells.Item(67, 5).Value = 1
$ws.Cells.Item(67, 6).Value = 'Occupied'
Set-TextCell $ws 68 1 61 '2026-02-01'
$ws.Cells.Item(68, 2).Value = '21:04:36'
$ws.Cells.Item(68, 3).Value = '21:00'
$ws.Cells.Item(68, 4).Value = 'Bedroom'
$ws.Cells.Item(68, 5).Value = 20
$ws.Cells.Item(68, 6).Value = 'Occupied'
Set-TextCell $ws 69 1 61 '2026-02-01'
$ws.Cells.Item(69, 2).Value = '21:04:41'
$ws.Cells.Item(69, 3).Value = '21:00'
$ws.Cells.Item(69, 4).Value = 'Bedroom'
$ws.Cells.Item(69, 5).Value = 2
$ws.Cells.Item(69, 6).Value = 'Occupied'
Set-TextCell $ws 70 1 61 '2026-02-01'
$ws.Cells.Item(70, 2).Value = '21:04:41'
$ws.Cells.Item(70, 3).Value = '21:00'
$ws.Cells.Item(70, 4).Value = 'Bedroom'
$ws.Cells.Item(70, 5).Value = 5
$ws.Cells.Item(70, 6).Value = 'Occupied'
Set-TextCell $ws 71 1 61 '2026-02-01'
$ws.Cells.Item(71, 2).Value = '21:04:42'
$ws.Cells.Item(71, 3).Value = '21:00'
$ws.Cells.Item(71, 4).Value = 'Bedroom'
$ws.Cells.Item(71, 5).Value = 2
$ws.Cells.Item(71, 6).Value = 'Occupied'
Set-TextCell $ws 72 1 61 '2026-02-01'
$ws.Cells.Item(72, 2).Value = '21:04:45'
$ws.Cells.Item(72, 3).Value = '21:00'
$ws.Cells.Item(72, 4).Value = 'Bedroom'
$ws.Cells.Item(72, 5).Value = 1
$ws.Cells.Item(72, 6).Value = 'Occupied'
Set-TextCell $ws 73 1 61 '2026-02-01'
$ws.Cells.Item(73, 2).Value = '21:04:50'
$ws.Cells.Item(73, 3).Value = '21:00'
$ws.Cells.Item(73, 4).Value = 'Bedroom'
$ws.Cells.Item(73, 5).Value = 2
$ws.Cells.Item(73, 6).Value = 'Occupied'
Set-TextCell $ws 74 1 61 '2026-02-01'
$ws.Cells.Item(74, 2).Value = '21:04:58'
$ws.Cells.Item(74, 3).Value = '21:00'
$ws.Cells.Item(74, 4).Value = 'Bedroom'
$ws.Cells.Item(74, 5).Value = 1
$ws.Cells.Item(74, 6).Value = 'Occupied'
Set-TextCell $ws 75 1 61 '2026-02-01'
$ws.Cells.Item(75, 2).Value = '21:04:59'
$ws.Cells.Item(75, 3).Value = '21:00'
$ws.Cells.Item(75, 4).Value = 'Bedroom'
$ws.Cells.Item(75, 5).Value = 2
$ws.Cells.Item(75, 6).Value = 'Occupied'
Set-TextCell $ws 76 1 61 '2026-02-01'
$ws.Cells.Item(76, 2).Value = '21:05:07'
$ws.Cells.Item(76, 3).Value = '21:00'
$ws.Cells.Item(76, 4).Value = 'Bedroom'
$ws.Cells.Item(76, 5).Value = 1
$ws.Cells.Item(76, 6).Value = 'Occupied'
Set-TextCell $ws 77 1 61 '2026-02-01'
$ws.Cells.Item(77, 2).Value = '21:05:14'
$ws.Cells.Item(77, 3).Value = '21:00'
$ws.Cells.Item(77, 4).Value = 'Bedroom'
$ws.Cells.Item(77, 5).Value = 2
$ws.Cells.Item(77, 6).Value = 'Occupied'
Set-TextCell $ws 78 1 61 '2026-02-01'
$ws.Cells.Item(78, 2).Value = '21:05:17'
$ws.Cells.Item(78, 3).Value = '21:00'
$ws.Cells.Item(78, 4).Value = 'Bedroom'
$ws.Cells.Item(78, 5).Value = 6
$ws.Cells.Item(78, 6).Value = 'Occupied'
Set-TextCell $ws 79 1 61 '2026-02-01'
$ws.Cells.Item(79, 2).Value = '21:05:18'
$ws.Cells.Item(79, 3).Value = '21:00'
$ws.Cells.Item(79, 4).Value = 'Bedroom'
$ws.Cells.Item(79, 5).Value = 2
$ws.Cells.Item(79, 6).Value = 'Occupied'
Set-TextCell $ws 80 1 61 '2026-02-01'
$ws.Cells.Item(80, 2).Value = '21:05:19'
$ws.Cells.Item(80, 3).Value = '21:00'
$ws.Cells.Item(80, 4).Value = 'Bedroom'
$ws.Cells.Item(80, 5).Value = 3
$ws.Cells.Item(80, 6).Value = 'Occupied'
Set-TextCell $ws 81 1 61 '2026-02-01'
$ws.Cells.Item(81, 2).Value = '21:05:20'
$ws.Cells.Item(81, 3).Value = '21:00'
$ws.Cells.Item(81, 4).Value = 'Bedroom'
$ws.Cells.Item(81, 5).Value = 2
$ws.Cells.Item(81, 6).Value = 'Occupied'

# ---- Sheet: mmWave(HR) ----
$ws = $wb.Worksheets.Item("mmWave(HR)")
Set-TextCell $ws 63 1 62 '2026-02-01'
$ws.Cells.Item(63, 2).Value = '21:04:23'
$ws.Cells.Item(63, 3).Value = '21:00'
$ws.Cells.Item(63, 4).Value = 'Bedroom'
$ws.Cells.Item(63, 5).Value = 57
$ws.Cells.Item(63, 6).Value = 'Occupied'
Set-TextCell $ws 64 1 62 '2026-02-01'
$ws.Cells.Item(64, 2).Value = '21:04:24'
$ws.Cells.Item(64, 3).Value = '21:00'
$ws.Cells.Item(64, 4).Value = 'Bedroom'
$ws.Cells.Item(64, 5).Value = 55
$ws.Cells.Item(64, 6).Value = 'Occupied'
Set-TextCell $ws 65 1 62 '2026-02-01'
$ws.Cells.Item(65, 2).Value = '21:04:24'
$ws.Cells.Item(65, 3).Value = '21:00'
$ws.Cells.Item(65, 4).Value = 'Bedroom'
$ws.Cells.Item(65, 5).Value = 50
$ws.Cells.Item(65, 6).Value = 'Occupied'
Set-TextCell $ws 66 1 62 '2026-02-01'
$ws.Cells.Item(66, 2).Value = '21:04:25'
$ws.Cells.Item(66, 3).Value = '21:00'
$ws.Cells.Item(66, 4).Value = 'Bedroom'
$ws.Cells.Item(66, 5).Value = 52
$ws.Cells.Item(66, 6).Value = 'Occupied'
Set-TextCell $ws 67 1 62 '2026-02-01'
$ws.Cells.Item(67, 2).Value = '21:04:26'
$ws.Cells.Item(67, 3).Value = '21:00'
$ws.Cells.Item(67, 4).Value = 'Bedroom'
$ws.Cells.Item(67, 5).Value = 50
$ws.Cells.Item(67, 6).Value = 'Occupied'
Set-TextCell $ws 68 1 62 '2026-02-01'
$ws.Cells.Item(68, 2).Value = '21:04:27'
$ws.Cells.Item(68, 3).Value = '21:00'
$ws.Cells.Item(68, 4).Value = 'Bedroom'
$ws.Cells.Item(68, 5).Value = 49
$ws.Cells.Item(68, 6).Value = 'Occupied'
Set-TextCell $ws 69 1 62 '2026-02-01'
$ws.Cells.Item(69, 2).Value = '21:04:36'
$ws.Cells.Item(69, 3).Value = '21:00'
$ws.Cells.Item(69, 4).Value = 'Bedroom'
$ws.Cells.Item(69, 5).Value = 68
$ws.Cells.Item(69, 6).Value = 'Occupied'
Set-TextCell $ws 70 1 62 '2026-02-01'
$ws.Cells.Item(70, 2).Value = '21:04:40'
$ws.Cells.Item(70, 3).Value = '21:00'
$ws.Cells.Item(70, 4).Value = 'Bedroom'
$ws.Cells.Item(70, 5).Value = 50
$ws.Cells.Item(70, 6).Value = 'Occupied'
Set-TextCell $ws 71 1 62 '2026-02-01'
$ws.Cells.Item(71, 2).Value = '21:04:41'
$ws.Cells.Item(71, 3).Value = '21:00'
$ws.Cells.Item(71, 4).Value = 'Bedroom'
$ws.Cells.Item(71, 5).Value = 53
$ws.Cells.Item(71, 6).Value = 'Occupied'
Set-TextCell $ws 72 1 62 '2026-02-01'
$ws.Cells.Item(72, 2).Value = '21:04:42'
$ws.Cells.Item(72, 3).Value = '21:00'
$ws.Cells.Item(72, 4).Value = 'Bedroom'
$ws.Cells.Item(72, 5).Value = 50
$ws.Cells.Item(72, 6).Value = 'Occupied'
Set-TextCell $ws 73 1 62 '2026-02-01'
$ws.Cells.Item(73, 2).Value = '21:04:45'
$ws.Cells.Item(73, 3).Value = '21:00'
$ws.Cells.Item(73, 4).Value = 'Bedroom'
$ws.Cells.Item(73, 5).Value = 49
$ws.Cells.Item(73, 6).Value = 'Occupied'
Set-TextCell $ws 74 1 62 '2026-02-01'
$ws.Cells.Item(74, 2).Value = '21:04:49'
$ws.Cells.Item(74, 3).Value = '21:00'
$ws.Cells.Item(74, 4).Value = 'Bedroom'
$ws.Cells.Item(74, 5).Value = 50
$ws.Cells.Item(74, 6).Value = 'Occupied'
Set-TextCell $ws 75 1 62 '2026-02-01'
$ws.Cells.Item(75, 2).Value = '21:04:57'
$ws.Cells.Item(75, 3).Value = '21:00'
$ws.Cells.Item(75, 4).Value = 'Bedroom'
$ws.Cells.Item(75, 5).Value = 49
$ws.Cells.Item(75, 6).Value = 'Occupied'
Set-TextCell $ws 76 1 62 '2026-02-01'
$ws.Cells.Item(76, 2).Value = '21:04:58'
$ws.Cells.Item(76, 3).Value = '21:00'
$ws.Cells.Item(76, 4).Value = 'Bedroom'
$ws.Cells.Item(76, 5).Value = 50
$ws.Cells.Item(76, 6).Value = 'Occupied'
Set-TextCell $ws 77 1 62 '2026-02-01'
$ws.Cells.Item(77, 2).Value = '21:05:06'
$ws.Cells.Item(77, 3).Value = '21:00'
$ws.Cells.Item(77, 4).Value = 'Bedroom'
$ws.Cells.Item(77, 5).Value = 49
$ws.Cells.Item(77, 6).Value = 'Occupied'
Set-TextCell $ws 78 1 62 '2026-02-01'
$ws.Cells.Item(78, 2).Value = '21:05:14'
$ws.Cells.Item(78, 3).Value = '21:00'
$ws.Cells.Item(78, 4).Value = 'Bedroom'
$ws.Cells.Item(78, 5).Value = 50
$ws.Cells.Item(78, 6).Value = 'Occupied'
Set-TextCell $ws 79 1 62 '2026-02-01'
$ws.Cells.Item(79, 2).Value = '21:05:17'
$ws.Cells.Item(79, 3).Value = '21:00'
$ws.Cells.Item(79, 4).Value = 'Bedroom'
$ws.Cells.Item(79, 5).Value = 54
$ws.Cells.Item(79, 6).Value = 'Occupied'
Set-TextCell $ws 80 1 62 '2026-02-01'
$ws.Cells.Item(80, 2).Value = '21:05:17'
$ws.Cells.Item(80, 3).Value = '21:00'
$ws.Cells.Item(80, 4).Value = 'Bedroom'
$ws.Cells.Item(80, 5).Value = 50
$ws.Cells.Item(80, 6).Value = 'Occupied'
Set-TextCell $ws 81 1 62 '2026-02-01'
$ws.Cells.Item(81, 2).Value = '21:05:18'
$ws.Cells.Item(81, 3).Value = '21:00'
$ws.Cells.Item(81, 4).Value = 'Bedroom'
$ws.Cells.Item(81, 5).Value = 51
$ws.Cells.Item(81, 6).Value = 'Occupied'
Set-TextCell $ws 82 1 62 '2026-02-01'
$ws.Cells.Item(82, 2).Value = '21:05:20'
$ws.Cells.Item(82, 3).Value = '21:00'
$ws.Cells.Item(82, 4).Value = 'Bedroom'
$ws.Cells.Item(82, 5).Value = 50
$ws.Cells.Item(82, 6).Value = 'Occupied'
Set-TextCell $ws 83 1 62 '2026-02-01'
$ws.Cells.Item(83, 2).Value = '21:05:23'
$ws.Cells.Item(83, 3).Value = '21:00'
$ws.Cells.Item(83, 4).Value = 'Bedroom'
$ws.Cells.Item(83, 5).Value = 60
$ws.Cells.Item(83, 6).Value = 'Occupied'

# ---- Sheet: mmWave(InBed) ----
$ws = $wb.Worksheets.Item("mmWave(InBed)")
Set-TextCell $ws 66 1 65 '2026-02-01'
$ws.Cells.Item(66, 2).Value = '21:04:23'
$ws.Cells.Item(66, 3).Value = '21:00'
$ws.Cells.Item(66, 4).Value = 'Bedroom'
$ws.Cells.Item(66, 5).Value = 'In Bed'
$ws.Cells.Item(66, 6).Value = 'Occupied'
Set-TextCell $ws 67 1 65 '2026-02-01'
$ws.Cells.Item(67, 2).Value = '21:04:23'
$ws.Cells.Item(67, 3).Value = '21:00'
$ws.Cells.Item(67, 4).Value = 'Bedroom'
$ws.Cells.Item(67, 5).Value = 'In Bed'
$ws.Cells.Item(67, 6).Value = 'Occupied'
Set-TextCell $ws 68 1 65 '2026-02-01'
$ws.Cells.Item(68, 2).Value = '21:04:24'
$ws.Cells.Item(68, 3).Value = '21:00'
$ws.Cells.Item(68, 4).Value = 'Bedroom'
$ws.Cells.Item(68, 5).Value = 'In Bed'
$ws.Cells.Item(68, 6).Value = 'Occupied'
Set-TextCell $ws 69 1 65 '2026-02-01'
$ws.Cells.Item(69, 2).Value = '21:04:25'
$ws.Cells.Item(69, 3).Value = '21:00'
$ws.Cells.Item(69, 4).Value = 'Bedroom'
$ws.Cells.Item(69, 5).Value = 'In Bed'
$ws.Cells.Item(69, 6).Value = 'Occupied'
Set-TextCell $ws 70 1 65 '2026-02-01'
$ws.Cells.Item(70, 2).Value = '21:04:26'
$ws.Cells.Item(70, 3).Value = '21:00'
$ws.Cells.Item(70, 4).Value = 'Bedroom'
$ws.Cells.Item(70, 5).Value = 'In Bed'
$ws.Cells.Item(70, 6).Value = 'Occupied'
Set-TextCell $ws 71 1 65 '2026-02-01'
$ws.Cells.Item(71, 2).Value = '21:04:27'
$ws.Cells.Item(71, 3).Value = '21:00'
$ws.Cells.Item(71, 4).Value = 'Bedroom'
$ws.Cells.Item(71, 5).Value = 'In Bed'
$ws.Cells.Item(71, 6).Value = 'Occupied'
Set-TextCell $ws 72 1 65 '2026-02-01'
$ws.Cells.Item(72, 2).Value = '21:04:35'
$ws.Cells.Item(72, 3).Value = '21:00'
$ws.Cells.Item(72, 4).Value = 'Bedroom'
$ws.Cells.Item(72, 5).Value = 'In Bed'
$ws.Cells.Item(72, 6).Value = 'Occupied'
Set-TextCell $ws 73 1 65 '2026-02-01'
$ws.Cells.Item(73, 2).Value = '21:04:40'
$ws.Cells.Item(73, 3).Value = '21:00'
$ws.Cells.Item(73, 4).Value = 'Bedroom'
$ws.Cells.Item(73, 5).Value = 'In Bed'
$ws.Cells.Item(73, 6).Value = 'Occupied'
Set-TextCell $ws 74 1 65 '2026-02-01'
$ws.Cells.Item(74, 2).Value = '21:04:41'
$ws.Cells.Item(74, 3).Value = '21:00'
$ws.Cells.Item(74, 4).Value = 'Bedroom'
$ws.Cells.Item(74, 5).Value = 'In Bed'
$ws.Cells.Item(74, 6).Value = 'Occupied'
Set-TextCell $ws 75 1 65 '2026-02-01'
$ws.Cells.Item(75, 2).Value = '21:04:42'
$ws.Cells.Item(75, 3).Value = '21:00'
$ws.Cells.Item(75, 4).Value = 'Bedroom'
$ws.Cells.Item(75, 5).Value = 'In Bed'
$ws.Cells.Item(75, 6).Value = 'Occupied'
Set-TextCell $ws 76 1 65 '2026-02-01'
$ws.Cells.Item(76, 2).Value = '21:04:44'
$ws.Cells.Item(76, 3).Value = '21:00'
$ws.Cells.Item(76, 4).Value = 'Bedroom'
$ws.Cells.Item(76, 5).Value = 'In Bed'
$ws.Cells.Item(76, 6).Value = 'Occupied'
Set-TextCell $ws 77 1 65 '2026-02-01'
$ws.Cells.Item(77, 2).Value = '21:04:49'
$ws.Cells.Item(77, 3).Value = '21:00'
$ws.Cells.Item(77, 4).Value = 'Bedroom'
$ws.Cells.Item(77, 5).Value = 'In Bed'
$ws.Cells.Item(77, 6).Value = 'Occupied'
Set-TextCell $ws 78 1 65 '2026-02-01'
$ws.Cells.Item(78, 2).Value = '21:04:57'
$ws.Cells.Item(78, 3).Value = '21:00'
$ws.Cells.Item(78, 4).Value = 'Bedroom'
$ws.Cells.Item(78, 5).Value = 'In Bed'
$ws.Cells.Item(78, 6).Value = 'Occupied'
Set-TextCell $ws 79 1 65 '2026-02-01'
$ws.Cells.Item(79, 2).Value = '21:04:58'
$ws.Cells.Item(79, 3).Value = '21:00'
$ws.Cells.Item(79, 4).Value = 'Bedroom'
$ws.Cells.Item(79, 5).Value = 'In Bed'
$ws.Cells.Item(79, 6).Value = 'Occupied'
Set-TextCell $ws 80 1 65 '2026-02-01'
$ws.Cells.Item(80, 2).Value = '21:05:06'
$ws.Cells.Item(80, 3).Value = '21:00'
$ws.Cells.Item(80, 4).Value = 'Bedroom'
$ws.Cells.Item(80, 5).Value = 'In Bed'
$ws.Cells.Item(80, 6).Value = 'Occupied'
Set-TextCell $ws 81 1 65 '2026-02-01'
$ws.Cells.Item(81, 2).Value = '21:05:13'
$ws.Cells.Item(81, 3).Value = '21:00'
$ws.Cells.Item(81, 4).Value = 'Bedroom'
$ws.Cells.Item(81, 5).Value = 'In Bed'
$ws.Cells.Item(81, 6).Value = 'Occupied'
Set-TextCell $ws 82 1 65 '2026-02-01'
$ws.Cells.Item(82, 2).Value = '21:05:16'
$ws.Cells.Item(82, 3).Value = '21:00'
$ws.Cells.Item(82, 4).Value = 'Bedroom'
$ws.Cells.Item(82, 5).Value = 'In Bed'
$ws.Cells.Item(82, 6).Value = 'Occupied'
Set-TextCell $ws 83 1 65 '2026-02-01'
$ws.Cells.Item(83, 2).Value = '21:05:17'
$ws.Cells.Item(83, 3).Value = '21:00'
$ws.Cells.Item(83, 4).Value = 'Bedroom'
$ws.Cells.Item(83, 5).Value = 'In Bed'
$ws.Cells.Item(83, 6).Value = 'Occupied'
Set-TextCell $ws 84 1 65 '2026-02-01'
$ws.Cells.Item(84, 2).Value = '21:05:18'
$ws.Cells.Item(84, 3).Value = '21:00'
$ws.Cells.Item(84, 4).Value = 'Bedroom'
$ws.Cells.Item(84, 5).Value = 'In Bed'
$ws.Cells.Item(84, 6).Value = 'Occupied'
Set-TextCell $ws 85 1 65 '2026-02-01'
$ws.Cells.Item(85, 2).Value = '21:05:20'
$ws.Cells.Item(85, 3).Value = '21:00'
$ws.Cells.Item(85, 4).Value = 'Bedroom'
$ws.Cells.Item(85, 5).Value = 'In Bed'
$ws.Cells.Item(85, 6).Value = 'Occupied'
Set-TextCell $ws 86 1 65 '2026-02-01'
$ws.Cells.Item(86, 2).Value = '21:05:22'
$ws.Cells.Item(86, 3).Value = '21:00'
$ws.Cells.Item(86, 4).Value = 'Bedroom'
$ws.Cells.Item(86, 5).Value = 'In Bed'
$ws.Cells.Item(86, 6).Value = 'Occupied'

